# Insert a new weekly price record before the existing row 108
# (shifts every subsequent "Perejil" record down by one row) and
# populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(108).Insert()

$ws.Range("A108").Value = 8
$ws.Range("B108").Value = "Terminal La Palmera de La Serena"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44824
$ws.Range("E108").Value = 4
$ws.Range("F108").Value = 100112044
$ws.Range("G108").Value = "Perejil"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 2800
$ws.Range("K108").Value = 2000
$ws.Range("L108").Value = 2500
$ws.Range("M108").Value = 2250
$ws.Range("N108").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O108").Value = "Provincia del Elquí"
$ws.Range("P108").Value = 1500
$ws.Range("Q108").Value = 1.5
$ws.Range("R108").Value = "Hortaliza"
